# Updated cryptos list values (price + 1h volume change) per target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price values that look like plain numbers must be forced to Text
# format first so Excel stores them verbatim (e.g. keeps trailing "315.80"
# instead of silently normalising to 315.8). Values with two dots (e.g.
# "28.513.35") are never auto-parsed as numbers, so they do not need this.

$ws.Range("D2").Value = '28.513.35'
$ws.Range("E2").Value = '  +0.80%  '

$ws.Range("D3").Value = '1.875.53'
$ws.Range("E3").Value = '  +0.88%  '

$ws.Range("E4").Value = '  -0.82%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.80'
$ws.Range("E5").Value = '  +0.49%  '

$ws.Range("E6").Value = '  -0.48%  '

$ws.Range("E7").Value = '  -0.32%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3906'
$ws.Range("E8").Value = '  -0.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08368'
$ws.Range("E9").Value = '  +1.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.104'
$ws.Range("E10").Value = '  -0.72%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.58'
$ws.Range("E11").Value = '  -0.60%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.229'
$ws.Range("E12").Value = '  +0.31%  '

$ws.Range("D13").Value = '1.874.03'
$ws.Range("E13").Value = '  +0.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.38'
$ws.Range("E14").Value = '  +0.65%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.268'
$ws.Range("E15").Value = '  +1.04%  '

$ws.Range("E16").Value = '  -0.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001104'
$ws.Range("E17").Value = '  +0.51%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.27'
$ws.Range("E18").Value = '  +0.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06728'
$ws.Range("E19").Value = '  +0.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.69'
$ws.Range("E20").Value = '  +0.73%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.008'
$ws.Range("E21").Value = '  -0.48%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.917'
$ws.Range("E22").Value = '  -0.31%  '

$ws.Range("D23").Value = '28.529.60'
$ws.Range("E23").Value = '  +0.74%  '

$ws.Range("E24").Value = '  +0.52%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.226'
$ws.Range("E25").Value = '  -1.33%  '

$ws.Range("D26").Value = '2.088.06'
$ws.Range("E26").Value = '  +0.87%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.96'
$ws.Range("E27").Value = '  +1.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.61'
$ws.Range("E28").Value = '  +0.24%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.387'
$ws.Range("E29").Value = '  -0.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.59'
$ws.Range("E30").Value = '  -0.72%  '

$ws.Range("E31").Value = '  -0.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.044'
$ws.Range("E32").Value = '  +1.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.767'
$ws.Range("E33").Value = '  -0.37%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.617'
$ws.Range("E34").Value = '  -0.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02458'
$ws.Range("E35").Value = '  +1.30%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06527'
$ws.Range("E36").Value = '  +1.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2161'
$ws.Range("E37").Value = '  -0.21%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.811'
$ws.Range("E38").Value = '  -3.47%  '

$ws.Range("E39").Value = '  +2.27%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.191'
$ws.Range("E40").Value = '  +0.88%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.239'
$ws.Range("E41").Value = '  -1.18%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6389'
$ws.Range("E42").Value = '  -0.27%  '

$ws.Range("E43").Value = '  -0.12%  '

$ws.Range("E44").Value = '  -0.49%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6001'
$ws.Range("E45").Value = '  +0.35%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.02'
$ws.Range("E46").Value = '  +1.11%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.683'
$ws.Range("E47").Value = '  -0.15%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.006'
$ws.Range("E48").Value = '  +1.48%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.215'
$ws.Range("E49").Value = '  +1.29%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '121.87'
$ws.Range("E50").Value = '  +1.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.144'
$ws.Range("E51").Value = '  -10.75%  '
